$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + volume-change columns) to match the latest
# coinranking.com snapshot.
#
# NOTE: several Price (column D) values look like plain numbers (e.g.
# "0.999", "8.04"). The sheet stores them as text, so a direct
# `Range.Value = "0.999"` assignment must be avoided: Excel's COM layer
# auto-coerces a numeric-looking string into a real number (and picking a
# Text NumberFormat afterwards leaves a stray style/quote-prefix behind).
# Routing the literal through a text-producing formula and then pasting
# back as values-only keeps the cell a plain string with no style change,
# which matches how the original text cells were authored.

$ws.Range("D2").Value = "62.638.92"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.448.37"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").Formula = '="0.999"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Formula = '="578.33"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Formula = '="148.56"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Formula = '="0.485"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Formula = '="8.04"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +5.30%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").Value = "4.039.39"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Formula = '="28.32"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -5.61%  "
$ws.Range("D15").Value = "3.447.46"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "62.681.62"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Formula = '="6.41"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="14.59"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Formula = '="387.39"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Formula = '="0.567"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "3.584.54"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Formula = '="0.183"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Formula = '="8.00"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").Formula = '="1.62"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="31.89"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Formula = '="6.94"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").Formula = '="169.24"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "3.483.33"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Formula = '="0.0775"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Formula = '="42.71"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Formula = '="4.39"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -2.59%  "
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "2.568.84"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("E51").Value = "  -0.01%  "
$excel.CutCopyMode = 0
